$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the paragraph "I could've spent hours writing permutation tests
#    ... in a BDD environment for the specification." entirely (including its
#    paragraph mark), so the "To test customer prices..." paragraph moves up
#    to directly follow the "...opening new routes/editing routes." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I could*ve spent hours writing permutation tests*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Delete the paragraph "The method inside KPSServer getTransportMap()
#    .calculateRoute(Mail) ... could only return one or the other." entirely.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The method inside KPSServer*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Remove the stray <w:lastRenderedPageBreak/> rendering hint that sits in
#    front of "The reason for this is directly after ...". Re-writing the
#    run's text via Find/Replace drops the stale rendering artifact while
#    keeping the visible text identical.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Execute(
    "The reason for this is directly after",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "The reason for this is directly after", 2)

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the end of the "...opening new
#    routes/editing routes." paragraph to the middle of "should've" (between
#    "shou" and "ld've") in the "I was never sure ..." paragraph. Re-adding a
#    bookmark with the same reserved name relocates it (Word keeps a single
#    "_GoBack" bookmark) and splits the run precisely at the given point.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("I was never sure if I shou")
$splitPoint = $rng2.End
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
